$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the odds in row 3 (M3: 1.13 -> 1.11, N3: 6 -> 6.5)
$ws.Range("M3").Value = 1.11
$ws.Range("N3").Value = 6.5

# Delete row 4 (nyeDPXV6 / Botafogo SP vs Avai) entirely, shifting rows 5 and 6 up
$ws.Rows.Item(4).Delete()
